$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Angel Warrior', ['Token Creature — Angel Warrior', 'Flying, vigilance', '4/4'])"
$ws.Range("A3").Value = "('Bear', ['Token Creature — Bear', '2/2'])"
$ws.Range("A4").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A5").Value = "('Cat', ['Token Creature — Cat', '2/2'])"
$ws.Range("A6").Value = "('Demon Berserker', ['Token Creature — Demon Berserker', 'Menace (This creature can’t be blocked except by two or more creatures.)', '2/3'])"
$ws.Range("A7").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '5/5'])"
$ws.Range("A8").Value = "('Dwarf Berserker', ['Token Creature — Dwarf Berserker', '2/1'])"
$ws.Range("A9").Value = "('Elf Warrior', ['Token Creature — Elf Warrior', '1/1'])"
$ws.Range("A10").Value = "('Foretell', ['Card', '(After you foretell a card, you can place the exiled card here. You may cast it on a later turn for its foretell cost.)'])"
$ws.Range("A11").Value = "('Giant Wizard', ['Token Creature — Giant Wizard', '4/4'])"
$ws.Range("A12").Value = "('Human Warrior', ['Token Creature — Human Warrior', '1/1'])"
$ws.Range("A13").Value = "('Icy Manalith', ['Token Snow Artifact', '{T}: Add one mana of any color.'])"
$ws.Range("A14").Value = "('Kaya the Inexorable Emblem', ['Emblem', 'At the beginning of your upkeep, you may cast a legendary spell from your hand, from your graveyard, or from among cards you own in exile without paying its mana cost.'])"
$ws.Range("A15").Value = "(`"Koma's Coil`", ['Token Creature — Serpent', '3/3'])"
$ws.Range("A16").Value = "('Replicated Ring', ['Token Snow Artifact', '{T}: Add one mana of any color.'])"
$ws.Range("A17").Value = "('Shapeshifter', ['Token Creature — Shapeshifter', 'Changeling (This token is every creature type.)', '2/2'])"
$ws.Range("A18").Value = "('Shard', ['Token Enchantment — Shard', '{2}, Sacrifice this enchantment: Scry 1, then draw a card.'])"
$ws.Range("A19").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A20").Value = "('Tibalt, Cosmic Impostor Emblem', ['Emblem', 'You may play cards exiled with Tibalt, Cosmic Impostor, and you may spend mana as though it were mana of any color to cast those spells.'])"
$ws.Range("A21").Value = "('Treasure', ['Token Artifact — Treasure', '{T}, Sacrifice this artifact: Add one mana of any color.'])"
$ws.Range("A22").Value = "('Troll Warrior', ['Token Creature — Troll Warrior', 'Trample', '4/4'])"
$ws.Range("A23").Value = "('Tyvar Kell Emblem', ['Emblem', 'Whenever you cast an Elf spell, it gains haste until end of turn and you draw two cards.'])"
$ws.Range("A24").Value = "('Zombie Berserker', ['Token Creature — Zombie Berserker', '2/2'])"
